# Auto-generated Excel COM-interop script
# Commit: Add data for 2024-12-24
# Applies the 2024-12-24 daily crime-count increments across the
# "Citywide Totals", "By Neighborhood", and per-neighborhood sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("K2").Value = 7748
$ws.Range("K3").Value = 8002
$ws.Range("J4").Value = 1850
$ws.Range("K4").Value = 1685
$ws.Range("K5").Value = 574
$ws.Range("K6").Value = 8935
$ws.Range("J7").Value = 29320
$ws.Range("K7").Value = 26944

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("K3").Value = 91
$ws.Range("K6").Value = 141
$ws.Range("K7").Value = 342

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("K2").Value = 491
$ws.Range("K3").Value = 532
$ws.Range("K6").Value = 593
$ws.Range("K7").Value = 1763

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("K5").Value = 15
$ws.Range("K7").Value = 570

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("K6").Value = 107
$ws.Range("K7").Value = 448

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("K2").Value = 252
$ws.Range("K3").Value = 295
$ws.Range("K7").Value = 887

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("K2").Value = 230
$ws.Range("K7").Value = 804
$ws.Range("K8").Value = 1763
$ws.Range("K9").Value = 128
$ws.Range("K11").Value = 473
$ws.Range("K18").Value = 183
$ws.Range("K19").Value = 775
$ws.Range("K20").Value = 664
$ws.Range("K21").Value = 92
$ws.Range("K25").Value = 127
$ws.Range("K29").Value = 1489
$ws.Range("K31").Value = 321
$ws.Range("K36").Value = 349
$ws.Range("K37").Value = 887
$ws.Range("K42").Value = 1000
$ws.Range("K47").Value = 183
$ws.Range("K48").Value = 336
$ws.Range("K49").Value = 150
$ws.Range("K51").Value = 348
$ws.Range("K53").Value = 342
$ws.Range("K55").Value = 295
$ws.Range("J63").Value = 197
$ws.Range("K63").Value = 78
$ws.Range("K64").Value = 162
$ws.Range("K67").Value = 1045
$ws.Range("K71").Value = 82
$ws.Range("K77").Value = 177
$ws.Range("K79").Value = 659
$ws.Range("K80").Value = 103
$ws.Range("K83").Value = 570
$ws.Range("K85").Value = 1240
$ws.Range("K86").Value = 163
$ws.Range("K89").Value = 405
$ws.Range("K91").Value = 322
$ws.Range("K93").Value = 108
$ws.Range("K95").Value = 448
$ws.Range("K96").Value = 288
$ws.Range("K97").Value = 220
$ws.Range("K98").Value = 145
$ws.Range("J101").Value = 29320
$ws.Range("K101").Value = 26944

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("K6").Value = 127
$ws.Range("K7").Value = 321

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("K2").Value = 284
$ws.Range("K7").Value = 1045

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("K4").Value = 13
$ws.Range("K7").Value = 150

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("K2").Value = 418
$ws.Range("K3").Value = 524
$ws.Range("K6").Value = 442
$ws.Range("K7").Value = 1489

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("K2").Value = 51
$ws.Range("K7").Value = 336

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("K2").Value = 226
$ws.Range("K3").Value = 231
$ws.Range("K7").Value = 775

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("K2").Value = 267
$ws.Range("K6").Value = 383
$ws.Range("K7").Value = 1000

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("K6").Value = 109
$ws.Range("K7").Value = 295

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("K4").Value = 21
$ws.Range("K7").Value = 288

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("K2").Value = 86
$ws.Range("K3").Value = 147
$ws.Range("K7").Value = 322

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Range("K6").Value = 56
$ws.Range("K7").Value = 92

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("K2").Value = 216
$ws.Range("K7").Value = 659

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("K6").Value = 61
$ws.Range("K7").Value = 162

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("K6").Value = 192
$ws.Range("K7").Value = 664

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("K5").Value = 5
$ws.Range("K7").Value = 183

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("K2").Value = 131
$ws.Range("K7").Value = 349

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range("K3").Value = 25
$ws.Range("K6").Value = 42
$ws.Range("K7").Value = 108

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("K3").Value = 254
$ws.Range("K6").Value = 226
$ws.Range("K7").Value = 804

$ws = $wb.Worksheets.Item('East Side')
$ws.Range("K4").Value = 9
$ws.Range("K7").Value = 127

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("K6").Value = 56
$ws.Range("K7").Value = 183

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("K6").Value = 81
$ws.Range("K7").Value = 145

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("K2").Value = 163
$ws.Range("K7").Value = 473

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range("K2").Value = 42
$ws.Range("K7").Value = 128

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("K3").Value = 62
$ws.Range("K6").Value = 83

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("K2").Value = 70
$ws.Range("K7").Value = 230

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("K6").Value = 123
$ws.Range("K7").Value = 220

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("K3").Value = 125
$ws.Range("K4").Value = 47
$ws.Range("K7").Value = 405

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("K3").Value = 24
$ws.Range("K7").Value = 163

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("K2").Value = 94
$ws.Range("K3").Value = 96
$ws.Range("K4").Value = 38
$ws.Range("K7").Value = 348

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("K3").Value = 429
$ws.Range("K7").Value = 1240

$ws = $wb.Worksheets.Item('Oakland')
$ws.Range("K3").Value = 30
$ws.Range("K7").Value = 82

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("K6").Value = 27
$ws.Range("K7").Value = 177

$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Range("K6").Value = 52
$ws.Range("K7").Value = 103
